$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values,
# keeping them stored as text (matching the original inlineStr cells).

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "300.87"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "0.66%"
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "31.59"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "0.77%"
$cell.ClearFormats()

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.092"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "-1.33%"
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.07794"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "-2.93%"
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "2.238"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "-15.87%"
$cell.ClearFormats()

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "7.794"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "-0.66%"
$cell.ClearFormats()

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.830"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "0.15%"
$cell.ClearFormats()

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.9153"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "-0.32%"
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1752"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "1.02%"
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07541"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "4.56%"
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08976"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "8.53%"
$cell.ClearFormats()

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.03090"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "3.22%"
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.1003"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "0.63%"
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.001514"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "0.85%"
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.006036"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "-0.57%"
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.469"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "-0.87%"
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.248"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "-0.09%"
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.3293"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "0.29%"
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.1337"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "1.33%"
$cell.ClearFormats()

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.345"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "-6.34%"
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.1817"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "13.61%"
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.04586"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "0.01%"
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.001251"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "-0.85%"
$cell.ClearFormats()

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.004465"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "0.26%"
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0001249"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "5.85%"
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "-1.42%"
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01764"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "-4.43%"
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.04772"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "5.75%"
$cell.ClearFormats()

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.007627"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "8.43%"
$cell.ClearFormats()

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1358"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "1.06%"
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.002189"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "-2.30%"
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.01026"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "-1.75%"
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.00006205"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "-4.11%"
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "-0.07%"
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "28.79%"
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "-9.20%"
$cell.ClearFormats()

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.00002099"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "-0.07%"
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0001999"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "-0.07%"
$cell.ClearFormats()
